$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicators")

$ws.Range("B4").Value = 0.82407407407407407
$ws.Range("B5").Value = 0.8783783783783784
$ws.Range("B6").Value = 0.8666666666666667
$ws.Range("B7").Value = 0.87248322147651003
